$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (G17=38956)
$ws.Range("H17").Value = 594.8461
$ws.Range("J17").Value = 590.64
$ws.Range("L17").Value = 1771.92
$ws.Range("N17").Value = -2107.92

# Row 40 (G40=5505)
$ws.Range("H40").Value = 2070.238
$ws.Range("I40").Value = 2160.077
$ws.Range("J40").Value = 1924.25
$ws.Range("K40").Value = 2160.077
$ws.Range("L40").Value = 1924.25
$ws.Range("M40").Value = -1985.077
$ws.Range("N40").Value = -2274.25

# Row 55 (G55=5517)
$ws.Range("H55").Value = 284.2
$ws.Range("I55").Value = 212.57143
$ws.Range("J55").Value = 451.33334
$ws.Range("K55").Value = 212.57143
$ws.Range("L55").Value = 451.33334
$ws.Range("M55").Value = 1.428570000000008
$ws.Range("N55").Value = -879.33334

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G2=27713)
$ws.Range("H2").Value = 1197.8125
$ws.Range("I2").Value = 1225.1428
$ws.Range("J2").Value = 1006.5
$ws.Range("K2").Value = 1225.1428
$ws.Range("L2").Value = 1006.5
$ws.Range("M2").Value = -1112.1428
$ws.Range("N2").Value = -1232.5

# Row 32 (G32=44147)
$ws.Range("H32").Value = 4578.0464
$ws.Range("I32").Value = 4293.4053
$ws.Range("K32").Value = 4293.4053
$ws.Range("M32").Value = -4006.4053

# Row 45 (G45=27714)
$ws.Range("H45").Value = 1099.4166
$ws.Range("I45").Value = 631
$ws.Range("J45").Value = 1255.5555
$ws.Range("K45").Value = 631
$ws.Range("L45").Value = 1255.5555
$ws.Range("M45").Value = -254
$ws.Range("N45").Value = -2009.5555

# Row 46 (G46=3498)
$ws.Range("H46").Value = 3962.25
$ws.Range("I46").Value = 3499.5
$ws.Range("J46").Value = 4425
$ws.Range("K46").Value = 3499.5
$ws.Range("L46").Value = 4425
$ws.Range("M46").Value = -3180.5
$ws.Range("N46").Value = -5063

# Row 74 (G74=44000)
$ws.Range("H74").Value = 4448.5
$ws.Range("I74").Value = 900.5
$ws.Range("J74").Value = 9770.5
$ws.Range("K74").Value = 900.5
$ws.Range("L74").Value = 9770.5
$ws.Range("M74").Value = -26.5
$ws.Range("N74").Value = -11518.5

# Row 77 (G77=44000)
$ws.Range("H77").Value = 4448.5
$ws.Range("I77").Value = 900.5
$ws.Range("J77").Value = 9770.5
$ws.Range("K77").Value = 4502.5
$ws.Range("L77").Value = 48852.5
$ws.Range("M77").Value = -134.5
$ws.Range("N77").Value = -57588.5

# Row 116 (G116=27713)
$ws.Range("H116").Value = 1197.8125
$ws.Range("I116").Value = 1225.1428
$ws.Range("J116").Value = 1006.5
$ws.Range("K116").Value = 1225.1428
$ws.Range("L116").Value = 1006.5
$ws.Range("M116").Value = 1068.8572
$ws.Range("N116").Value = -5594.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G3=27713)
$ws.Range("H3").Value = 1197.8125
$ws.Range("I3").Value = 1225.1428
$ws.Range("J3").Value = 1006.5
$ws.Range("K3").Value = 1225.1428
$ws.Range("L3").Value = 1006.5
$ws.Range("M3").Value = -1111.1428
$ws.Range("N3").Value = -1234.5

# Row 22 (G22=5092)
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -627
$ws.Range("N22").Value = -646

# Row 99 (G99=19943)
$ws.Range("H99").Value = 967.1429000000001
$ws.Range("I99").Value = 961.6667
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 961.6667
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 536.3333
$ws.Range("N99").Value = -3996

# Row 134 (G134=43998)
$ws.Range("H134").Value = 15175706
$ws.Range("I134").Value = 15175706
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 45527118
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -45524583
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 102 (G102=19813)
$ws.Range("H102").Value = 2545.182
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 2699.7
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 8099.099999999999
$ws.Range("M102").Value = -566
$ws.Range("N102").Value = -12967.1

# Row 107 (G107=27838)
$ws.Range("H107").Value = 7865.125
$ws.Range("I107").Value = 362.875
$ws.Range("J107").Value = 15367.375
$ws.Range("K107").Value = 1088.625
$ws.Range("L107").Value = 46102.125
$ws.Range("M107").Value = 831.375
$ws.Range("N107").Value = -49942.125

# Row 108 (G108=27853)
$ws.Range("H108").Value = 400
$ws.Range("I108").Value = 400
$ws.Range("K108").Value = 1200
$ws.Range("M108").Value = 1680

# Row 110 (G110=27857)
$ws.Range("H110").Value = 4885.706
$ws.Range("I110").Value = 2523.4
$ws.Range("J110").Value = 5870
$ws.Range("K110").Value = 7570.200000000001
$ws.Range("L110").Value = 17610
$ws.Range("M110").Value = -3480.200000000001
$ws.Range("N110").Value = -25790

# Row 111 (G111=27856)
$ws.Range("H111").Value = 2439
$ws.Range("I111").Value = 1208.5
$ws.Range("J111").Value = 4900
$ws.Range("K111").Value = 3625.5
$ws.Range("L111").Value = 14700
$ws.Range("M111").Value = -558.5
$ws.Range("N111").Value = -20834

# Row 114 (G114=27865)
$ws.Range("H114").Value = 436.18182
$ws.Range("I114").Value = 429.8
$ws.Range("K114").Value = 1289.4
$ws.Range("M114").Value = 1964.6

# Row 117 (G117=27870)
$ws.Range("H117").Value = 100029
$ws.Range("I117").Value = 100029
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 300087
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -296645
$ws.Range("N117").ClearContents()

# Row 123 (G123=36037)
$ws.Range("H123").Value = 10156.25
$ws.Range("I123").Value = 18738.334
$ws.Range("J123").Value = 5007
$ws.Range("K123").Value = 56215.00199999999
$ws.Range("L123").Value = 15021
$ws.Range("M123").Value = -53765.00199999999
$ws.Range("N123").Value = -19921

# Row 129 (G129=36054)
$ws.Range("H129").Value = 1383.0385
$ws.Range("I129").Value = 472
$ws.Range("J129").Value = 1599.9524
$ws.Range("K129").Value = 1416
$ws.Range("L129").Value = 4799.857199999999
$ws.Range("M129").Value = 3584
$ws.Range("N129").Value = -14799.8572

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (G126=36184)
$ws.Range("H126").Value = 1638.1082
$ws.Range("I126").Value = 1234.6
$ws.Range("J126").Value = 2112.8235
$ws.Range("K126").Value = 3703.8
$ws.Range("L126").Value = 6338.470499999999
$ws.Range("M126").Value = -1233.8
$ws.Range("N126").Value = -11278.4705

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G22=5277)
$ws.Range("H22").Value = 1259
$ws.Range("I22").Value = 1341.6
$ws.Range("J22").Value = 983.6667
$ws.Range("K22").Value = 1341.6
$ws.Range("L22").Value = 983.6667
$ws.Range("M22").Value = -1046.6
$ws.Range("N22").Value = -1573.6667

# Row 27 (G27=5277)
$ws.Range("H27").Value = 1259
$ws.Range("I27").Value = 1341.6
$ws.Range("J27").Value = 983.6667
$ws.Range("K27").Value = 1341.6
$ws.Range("L27").Value = 983.6667
$ws.Range("M27").Value = -1234.6
$ws.Range("N27").Value = -1197.6667

# Row 46 (G46=5282)
$ws.Range("H46").Value = 1057.1428
$ws.Range("I46").Value = 1062.5
$ws.Range("J46").Value = 1040
$ws.Range("K46").Value = 1062.5
$ws.Range("L46").Value = 1040
$ws.Range("M46").Value = -874.5
$ws.Range("N46").Value = -1416

# Row 55 (G55=5284)
$ws.Range("H55").Value = 1820.7142
$ws.Range("I55").Value = 3850.3333
$ws.Range("J55").Value = 298.5
$ws.Range("K55").Value = 3850.3333
$ws.Range("L55").Value = 298.5
$ws.Range("M55").Value = -3677.3333
$ws.Range("N55").Value = -644.5
